$d = $word.ActiveDocument

# --- Update the date line at the top of the document ---
$d.Content.Find.Execute("2024-04-29 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-30 Tuesday", 2)

# --- Update the two-digit multiplication problems in the table ---
# Cells are addressed by (row, column) rather than by searching for the
# old text, because some new values coincide with other old values
# elsewhere in the table (e.g. "65×74=" is both an old and a new value),
# which would make a blind find-and-replace ambiguous/unsafe.
$t = $d.Tables.Item(1)

$updates = @(
    @{r=1;  c=1; v="60×66="},
    @{r=1;  c=2; v="29×41="},
    @{r=1;  c=3; v="44×38="},
    @{r=1;  c=4; v="56×49="},
    @{r=1;  c=5; v="82×58="},

    @{r=5;  c=1; v="91×57="},
    @{r=5;  c=2; v="96×66="},
    @{r=5;  c=3; v="66×58="},
    @{r=5;  c=4; v="13×33="},
    @{r=5;  c=5; v="97×98="},

    @{r=10; c=1; v="37×28="},
    @{r=10; c=2; v="46×65="},
    @{r=10; c=3; v="43×21="},
    @{r=10; c=4; v="36×14="},
    @{r=10; c=5; v="54×55="},

    @{r=15; c=1; v="65×74="},
    @{r=15; c=2; v="74×13="},
    @{r=15; c=3; v="37×76="},
    @{r=15; c=4; v="69×24="},
    @{r=15; c=5; v="31×49="},

    @{r=20; c=1; v="63×70="},
    @{r=20; c=2; v="86×73="},
    @{r=20; c=3; v="46×37="},
    @{r=20; c=4; v="21×69="},
    @{r=20; c=5; v="85×24="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.r, $u.c)
    $cell.Range.Text = $u.v
}
